$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder name / card number
$ws.Range("C2").Value = "Hartmut"
# Card number is a 16-digit string that must stay text (plain .Value
# assignment would auto-convert it to a floating point number and lose
# precision). Build the text in a scratch cell via a formula, then paste
# only the computed value into B3 so the destination keeps its existing
# (General-formatted) cell style untouched.
$ws.Range("Z1").Formula = "=""2570314725427075"""
$ws.Range("Z1").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 07.11.2023"

# Row 6
$ws.Range("B6").Value = "11.11."
$ws.Range("C6").Value = "12.11."
$ws.Range("D6").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E6").Value = "54,38-"

# Row 7
$ws.Range("B7").Value = "12.11."
$ws.Range("C7").Value = "13.11."
$ws.Range("D7").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E7").Value = "25,20-"

# Row 8
$ws.Range("B8").Value = "13.11."
$ws.Range("C8").Value = "14.11."
$ws.Range("D8").Value = "EBAY MKTPLC EU VAKLKS"
$ws.Range("E8").Value = "204,45-"

# Row 9
$ws.Range("B9").Value = "16.11."
$ws.Range("C9").Value = "17.11."
$ws.Range("D9").Value = "BEITRAG Allianz SE K-72518086"
$ws.Range("E9").Value = "54,60-"

# Row 10
$ws.Range("B10").Value = "18.11."
$ws.Range("C10").Value = "19.11."
$ws.Range("D10").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 85718469"
$ws.Range("E10").Value = "84,46-"

# Row 11 was empty (style 8/8/8/12); now gets a new transaction and E11 adopts
# the right-aligned, non-wrapping style used by the other amount cells (E6:E10)
$ws.Range("B11").Value = "20.11."
$ws.Range("C11").Value = "21.11."
$ws.Range("D11").Value = "PAYPAL XANNKU"
$ws.Range("E11").Value = "43,91-"
$ws.Range("E11").HorizontalAlignment = -4152
$ws.Range("E11").VerticalAlignment = -4107
$ws.Range("E11").WrapText = $false

# Row 12: closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 25.11.2023"
$ws.Range("E12").Value = "467,00-"

# Row 13: next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 30.11.2023"
